$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4550261933023185
$ws.Range("C2").Value = 0.7785243494829454
$ws.Range("D2").Value = 0.7737919637881365
$ws.Range("E2").Value = 0.3901673830723808
$ws.Range("F2").Value = 0.5359328471787868
$ws.Range("G2").Value = 0.4929636663620077
$ws.Range("H2").Value = -168.6401250326047
$ws.Range("I2").Value = 347.2802500652093
$ws.Range("J2").Value = 364.6834446819193
$ws.Range("B3").Value = 0.431129791961647
$ws.Range("C3").Value = 0.8025829984853882
$ws.Range("D3").Value = 0.7974992988755698
$ws.Range("E3").Value = 0.3769325639589143
$ws.Range("F3").Value = 0.5699942986560207
$ws.Range("G3").Value = 0.5213144079378343
$ws.Range("H3").Value = -154.8533735312403
$ws.Range("I3").Value = 321.7067470624805
$ws.Range("J3").Value = 342.5905806025325
$ws.Range("B4").Value = 0.4069223977525986
$ws.Range("C4").Value = 0.8240156243527104
$ws.Range("D4").Value = 0.8187057509495594
$ws.Range("E4").Value = 0.3748672222923064
$ws.Range("F4").Value = 0.5776741876625082
$ws.Range("G4").Value = 0.5208226360016921
$ws.Range("H4").Value = -140.9889179471134
$ws.Range("I4").Value = 295.9778358942269
$ws.Range("J4").Value = 320.3423083576208
$ws.Range("B5").Value = 0.3863824836961171
$ws.Range("C5").Value = 0.8411709234032838
$ws.Range("D5").Value = 0.8356703493220122
$ws.Range("E5").Value = 0.3605917587162446
$ws.Range("F5").Value = 0.6132051280381245
$ws.Range("G5").Value = 0.5525314226323401
$ws.Range("H5").Value = -128.7852963864772
$ws.Range("I5").Value = 273.5705927729544
$ws.Range("J5").Value = 301.4157041596903
$ws.Range("B6").Value = 0.3658276497474408
$ws.Range("C6").Value = 0.8579678347371049
$ws.Range("D6").Value = 0.8524100543572525
$ws.Range("E6").Value = 0.355523965275849
$ws.Range("F6").Value = 0.6302295448495192
$ws.Range("G6").Value = 0.5636708629224326
$ws.Range("H6").Value = -115.5945903272773
$ws.Range("I6").Value = 249.1891806545546
$ws.Range("J6").Value = 280.5149309646325
$ws.Range("B7").Value = 0.3467238624949215
$ws.Range("C7").Value = 0.8708034013737245
$ws.Range("D7").Value = 0.8651616285079482
$ws.Range("E7").Value = 0.3488277462030486
$ws.Range("F7").Value = 0.6410037740252483
$ws.Range("G7").Value = 0.5677392381120336
$ws.Range("H7").Value = -104.208762653963
$ws.Range("I7").Value = 228.417525307926
$ws.Range("J7").Value = 263.2239145413459
$ws.Range("B8").Value = 0.3347388592432293
$ws.Range("C8").Value = 0.8803131575107868
$ws.Range("D8").Value = 0.8745387923029738
$ws.Range("E8").Value = 0.3356548522859122
$ws.Range("F8").Value = 0.6708047414951064
$ws.Range("G8").Value = 0.5953641614210683
$ws.Range("H8").Value = -94.9956602040921
$ws.Range("I8").Value = 211.9913204081842
$ws.Range("J8").Value = 250.2783485649461
$ws.Range("B9").Value = 0.3183415625138936
$ws.Range("C9").Value = 0.8915324249059252
$ws.Range("D9").Value = 0.8857984561784852
$ws.Range("E9").Value = 0.333200625905965
$ws.Range("F9").Value = 0.6785280148092026
$ws.Range("G9").Value = 0.596450061143467
$ws.Range("H9").Value = -83.17789812405405
$ws.Range("I9").Value = 190.3557962481081
$ws.Range("J9").Value = 232.123463328212